$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.211.93"
$ws.Range("E2").Value = "  -3.17%  "

$ws.Range("D3").Value = "3.541.42"
$ws.Range("E3").Value = "  -3.67%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.26%  "

$ws.Range("D7").Value = "3.540.58"
$ws.Range("E7").Value = "  -3.57%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -2.43%  "

$ws.Range("E10").Value = "  -2.01%  "

$ws.Range("E11").Value = "  -3.17%  "

$ws.Range("E12").Value = "  -3.61%  "

$ws.Range("E13").Value = "  -4.07%  "

$ws.Range("D14").Value = "4.142.67"
$ws.Range("E14").Value = "  -3.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").Value = "3.555.59"
$ws.Range("E16").Value = "  -3.21%  "

$ws.Range("D17").Value = "67.193.33"
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("E19").Value = "  -1.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "454.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.51%  "

$ws.Range("E22").Value = "  -4.87%  "

$ws.Range("E23").Value = "  -0.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").Value = "3.682.73"
$ws.Range("E25").Value = "  -3.64%  "

$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000125"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.70%  "

$ws.Range("E30").Value = "  -1.87%  "

$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("E33").Value = "  -3.29%  "

$ws.Range("E34").Value = "  -5.13%  "

$ws.Range("E35").Value = "  -3.94%  "

$ws.Range("E36").Value = "  -4.45%  "

$ws.Range("D37").Value = "3.539.60"
$ws.Range("E37").Value = "  -3.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.13%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.44%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.09%  "

$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0876"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.71%  "

$ws.Range("E45").Value = "  -3.50%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.16%  "

$ws.Range("E50").Value = "  -2.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.22%  "
